$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# Change 1 ("Political Ideology" paragraph):
# "...correlate highly with one another (r > .7) then we will average..."
# becomes
# "...correlate highly with one another then we will average..."
# (the three runs holding "(", "r", " > .7) " collapse away; no bookmark
# nearby so a single contiguous replace is safe)
# ----------------------------------------------------------------------
$d.Content.Find.Execute("one another (r > .7) then we will average", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "one another then we will average", 2)

# ----------------------------------------------------------------------
# Change 2 ("Political Identity" paragraph):
# "...correlate highly (r > .7) <bookmark>then we will average them to
#  create a 'General Political Identity' score..."
# becomes
# "...correlate highly <bookmark>then we will average them to create a
#  'General Political Identity' score..."
# The bookmark (_GoBack) sits right after "(r > .7) ", so only that exact
# substring is removed, leaving the bookmark untouched.
# ----------------------------------------------------------------------
$d.Content.Find.Execute("correlate highly (r > .7) ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "correlate highly ", 2)

# The word "Political " needs to move from the run that currently reads
# "Political Identity' score, ..." to the end of the preceding run that
# reads "then we will average them to create a 'General ". Net text is
# unchanged overall, only where the run boundary falls.
$d.Content.Find.Execute("a `u{2018}General ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "a `u{2018}General Political ", 2)
$d.Content.Find.Execute("Political Identity`u{2019} score, with larger values", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Identity`u{2019} score, with larger values", 2)
